# Updated cryptos list - apply price and volume(1h) changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.905.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.832.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6900"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.98%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07691"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3049"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.64%  "

$ws.Range("E10").Value = "  -4.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07805"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.842.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.078"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.84%  "

$ws.Range("E14").Value = "  -3.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6808"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.447"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.12%  "

$ws.Range("E17").Value = "  -1.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.915.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.077.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.95%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.469"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.62%  "

$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1469"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.799"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.42%  "

$ws.Range("E28").Value = "  -3.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.543"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.209"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.147"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.185"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05106"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7655"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.840"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.141"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.686"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.85%  "

$ws.Range("E38").Value = "  -1.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.221.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.699"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9368"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9994"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.695"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000123"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.49%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.570"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5172"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.978.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.747"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.09%  "

$ws.Range("E51").Value = "  -2.74%  "
